$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "C2" (test) is removed
$ws.Range("C2").ClearContents()

# New row 3: Other Grant / CS_PI / Nirjon, Shahriar
$ws.Range("A3").Value = "Other Grant"
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Value = "CS_PI"
$ws.Range("C3").Value = "Nirjon, Shahriar"

# Row 2: "D2" (othertest) becomes "Roberts, Lee"
$ws.Range("D2").Value = "Roberts, Lee"

# Update selection to reflect the last-edited cell
[void]$ws.Range("D3").Select()
